# Grace 4th updates on 0804
# Applies updated "Output per Capita" row values and adds four new
# indicator rows (Education Cost, Healthcare Cost, Infrastructure
# Investment, Reinvestment), each expressed as % of GDP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (Output per Capita) : refreshed values ----
$row2Vals = @{
    "B2" = 4.87161297425764
    "C2" = 4.84802685386431
    "D2" = 4.8166776801437
    "E2" = 4.77612623694607
    "F2" = 4.73146342589992
    "G2" = 4.64801859576522
    "H2" = 4.56133679684499
    "I2" = 4.47275457175577
    "J2" = 4.3839986224852
    "K2" = 4.2996133347857
    "L2" = 4.21665059439685
    "M2" = 4.1334575501158
    "N2" = 4.05088171609849
    "O2" = 3.971503937343
    "P2" = 4.04343857409873
    "Q2" = 4.03732907851332
    "R2" = 3.97131126076895
    "S2" = 3.92798487130393
    "T2" = 3.86520477232934
    "U2" = 3.81851897805827
    "V2" = 3.74743830870591
    "W2" = 3.67882983984387
    "X2" = 3.63438719770882
    "Y2" = 3.56474952669185
    "Z2" = 3.49483157874768
    "AA2" = 3.42673814500145
    "AB2" = 3.39738112578505
    "AC2" = 3.33500375374842
    "AD2" = 3.27329391068869
}

# ---- Row 3 : new "Education Cost (%GDP)" series ----
$ws.Range("A3").Value = "'Education Cost (%GDP)"
$row3Vals = @{
    "B3" = 0.0189617558665476
    "C3" = 0.0177902064190956
    "D3" = 0.0167830195234511
    "E3" = 0.015872836536675
    "F3" = 0.0150318033926546
    "G3" = 0.0150369984839259
    "H3" = 0.015101350375476
    "I3" = 0.0151123594381412
    "J3" = 0.015192786964971
    "K3" = 0.0152570064127073
    "L3" = 0.0153707650129577
    "M3" = 0.0155800452000004
    "N3" = 0.0158275300444287
    "O3" = 0.0161101275570302
    "P3" = 0.0158127184239433
    "Q3" = 0.0158287660115228
    "R3" = 0.0160889760022968
    "S3" = 0.0162718954880892
    "T3" = 0.0165521074591424
    "U3" = 0.0168424913206976
    "V3" = 0.0172020573600097
    "W3" = 0.0175114793230743
    "X3" = 0.0176934126422211
    "Y3" = 0.0180421950358619
    "Z3" = 0.0183375114320541
    "AA3" = 0.0185304714816922
    "AB3" = 0.0184317791343575
    "AC3" = 0.0184438436126483
    "AD3" = 0.0183971529375592
}

# ---- Row 4 : new "Healthcare Cost (%GDP)" series ----
$ws.Range("A4").Value = "'Healthcare Cost (%GDP)"
$row4Vals = @{
    "B4" = 0.0152547016486263
    "D4" = 0.0146297478383039
    "E4" = 0.0144230764615648
    "F4" = 0.014259149124476
    "G4" = 0.0144593856823192
    "H4" = 0.0147020255101967
    "I4" = 0.0149775674058206
    "J4" = 0.0152713648419669
    "K4" = 0.0155389084322633
    "L4" = 0.0158137372566904
    "M4" = 0.0161150823075959
    "N4" = 0.0164333591709379
    "O4" = 0.0167529110564366
    "P4" = 0.0164410593847658
    "Q4" = 0.0164501822204273
    "R4" = 0.0167058376102606
    "S4" = 0.0168712749105834
    "T4" = 0.0171257500343471
    "U4" = 0.017325283892178
    "V4" = 0.0176310170304514
    "W4" = 0.0179227211223843
    "X4" = 0.0181069993126021
    "Y4" = 0.0184380642220201
    "Z4" = 0.0187643809482546
    "AA4" = 0.0190658269806423
    "AB4" = 0.0191364207585739
    "AC4" = 0.0193806854784594
    "AD4" = 0.0196170082551715
}

# ---- Row 5 : new "Infrastructure Investment (%GDP)" series ----
$ws.Range("A5").Value = "'Infrastructure Investment (%GDP)"
$row5Vals = @{
    "O5" = 0.834752145560285
    "P5" = 0.842776687578206
    "R5" = 0.771256808311876
    "T5" = 0.91666246158243
    "U5" = 0.713296208766686
    "V5" = 0.817119797664799
    "X5" = 0.821686463339073
    "Y5" = 0.837438836692497
    "Z5" = 0.865027654798329
    "AB5" = 0.84407775257674
    "AC5" = 0.857313187405282
    "AD5" = 0.34549081912521
}

# ---- Row 6 : new "Reinvestment (%GDP)" series ----
$ws.Range("A6").Value = "'Reinvestment (%GDP)"
$row6Vals = @{
    "B6" = 0.123216457977028
    "C6" = 0.10679020687724
    "D6" = 0.120412767816642
    "E6" = 0.119295913450432
    "F6" = 0.118290952966872
    "G6" = 0.118496384617152
    "H6" = 0.118803376338161
    "I6" = 0.119089927298339
    "J6" = 0.119464152263366
    "K6" = 0.119795915303284
    "L6" = 0.120184502729984
    "M6" = 0.120695127970274
    "N6" = 0.121260889680635
    "O6" = 0.956615184173751
    "P6" = 0.964030465386915
    "Q6" = 0.121278948679807
    "R6" = 0.893051621924434
    "S6" = 0.122143170847555
    "T6" = 1.03934031907592
    "U6" = 0.836463983979561
    "V6" = 0.940952872055261
    "W6" = 0.124434200904525
    "X6" = 0.946486875293896
    "Y6" = 0.962919095950379
    "Z6" = 0.991129547178638
    "AA6" = 0.126596298939521
    "AB6" = 0.970645952469671
    "AC6" = 0.984137716496389
    "AD6" = 0.472504980317946
}

# ---- write all numeric values ----
$allVals = @($row2Vals, $row3Vals, $row4Vals, $row5Vals, $row6Vals)
foreach ($rowVals in $allVals) {
    foreach ($key in $rowVals.Keys) {
        $ws.Range($key).Value = $rowVals[$key]
    }
}

# ---- bump the sheet's outlineLevelRow (1 -> 5) ----
# Matches the source edit's sheetFormatPr outlineLevelRow change without
# leaving a stray per-row outline attribute: group a scratch row below the
# data, then delete it, so only the sheet-wide max outline level sticks.
$scratchRow = $ws.Rows.Item(100)
$scratchRow.OutlineLevel = 5
$scratchRow.Delete()
